$p = $ppt.ActivePresentation

# Delete slides 2, 3, and 4 (in reverse order to keep indices valid),
# leaving only the first slide ("ERP Post-Simulation Analysis").
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
